# Slide 9 ("Conclusion"), Content Placeholder 2:
# "As the model is exceeding the criteria of predicting 80% ..."
#   -> "As the model is meeting the criteria of predicting 80% ..."
#
# The word "exceeding" is replaced with "meeting" by re-typing the phrase
# "is exceeding " as "is meeting ", which is what causes PowerPoint to
# split the single run into three runs:
#   "As the model " | "is meeting " | "the criteria of predicting 80% ..."

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(9)
$shape = $s.Shapes.Item(4)
$tr = $shape.TextFrame.TextRange

$fullText = $tr.Text
$startPos = $fullText.IndexOf("is exceeding ") + 1

$target = $tr.Characters($startPos, 13)
$target.Text = "is meeting "
